# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45206 (2023-10-07) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$firstRow = 2
$lastRow = 459

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45206) {
        $cell.Value2 = 45208
    }
}
